$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(29).Insert()
$ws.Range("D29").Font.Name = "Arial"
$ws.Range("D29").Font.Size = 9
$ws.Range("D29").Value = "L78L05ABUTR "
$ws.Rows.Item(29).AutoFit()
Write-Host "RowHeight:" $ws.Rows.Item(29).RowHeight
